$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6471
$ws.Range("K3").Value = 6666
$ws.Range("B4").Value = 1416
$ws.Range("K4").Value = 1390
$ws.Range("K5").Value = 478
$ws.Range("K6").Value = 7341
$ws.Range("B7").Value = 18513
$ws.Range("K7").Value = 22346

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 446
$ws.Range("K7").Value = 1466

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 480

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 126
$ws.Range("K7").Value = 367

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 217
$ws.Range("K4").Value = 36
$ws.Range("K7").Value = 759

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 189
$ws.Range("K7").Value = 523

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 153
$ws.Range("K7").Value = 371

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 196
$ws.Range("K4").Value = 82
$ws.Range("K7").Value = 671
$ws.Range("K8").Value = 1466
$ws.Range("K11").Value = 412
$ws.Range("K14").Value = 113
$ws.Range("K15").Value = 230
$ws.Range("K16").Value = 57
$ws.Range("K19").Value = 654
$ws.Range("B22").Value = 52
$ws.Range("K26").Value = 31
$ws.Range("K29").Value = 1207
$ws.Range("K34").Value = 127
$ws.Range("K36").Value = 284
$ws.Range("K37").Value = 759
$ws.Range("K42").Value = 827
$ws.Range("K43").Value = 183
$ws.Range("K48").Value = 281
$ws.Range("K50").Value = 106
$ws.Range("K52").Value = 587
$ws.Range("K53").Value = 284
$ws.Range("K60").Value = 131
$ws.Range("K63").Value = 63
$ws.Range("K65").Value = 523
$ws.Range("K67").Value = 876
$ws.Range("K72").Value = 117
$ws.Range("K73").Value = 199
$ws.Range("K83").Value = 480
$ws.Range("K84").Value = 180
$ws.Range("K86").Value = 136
$ws.Range("K89").Value = 331
$ws.Range("K90").Value = 209
$ws.Range("K94").Value = 299
$ws.Range("K95").Value = 367
$ws.Range("K96").Value = 239
$ws.Range("K99").Value = 371
$ws.Range("B101").Value = 18513
$ws.Range("K101").Value = 22346

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 241
$ws.Range("K3").Value = 319
$ws.Range("K6").Value = 247
$ws.Range("K7").Value = 876

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 344
$ws.Range("K6").Value = 348
$ws.Range("K7").Value = 1207

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 193
$ws.Range("K3").Value = 197
$ws.Range("K7").Value = 654

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K3").Value = 26
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 223
$ws.Range("K3").Value = 251
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 307
$ws.Range("K7").Value = 827

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 75
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 111
$ws.Range("K7").Value = 284

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 220
$ws.Range("K3").Value = 221
$ws.Range("K5").Value = 27
$ws.Range("K7").Value = 671

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 50
$ws.Range("K7").Value = 127

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 299

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 412

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 66
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 199

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 196

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 95
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 23
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 209

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 49
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("B4").Value = 7
$ws.Range("B7").Value = 52

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 26
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 154
$ws.Range("K6").Value = 213
$ws.Range("K7").Value = 587

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 57
